# Generate Report for Handoff
# Rows 7, 8, 9, 10, 11, 13 (on "zh-cn" and "de-de") just finished a new
# handoff round: their Priority is now "ht" and the handoff/generate
# timestamps advance a little under a minute. This updates the Overview
# sheet's "Latest HO Xliff Generate Date" and each locale sheet's
# "Latest Handoff Datetime" + "Priority" columns for those rows.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 11, 13)

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-23 12:20:34"
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-23 12:20:29"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-23 12:20:34"
}
